# Generate Report for Archive
#
# 1) Update the status text from "Ready for handoff" to "In Translation"
#    everywhere it appears (Overview!E2:E3/F2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2) Narrow the status columns (Overview E:F, zh-cn C, de-de C) to their
#    new, narrower width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Replace the status text wherever it currently reads "Ready for handoff" ---
$sheets = @($overview, $zhcn, $dede)
foreach ($ws in $sheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE 1: reading the parameterized "Value" property requires calling
        # it like a method (Value()) in this interop - a bare ".Value" read
        # returns the property's method signature instead of the cell data.
        # NOTE 2: put the string literal on the LEFT of -eq so PowerShell
        # doesn't coerce it to the left operand's type (e.g. a Boolean cell
        # value would otherwise convert any non-empty string to $true and
        # falsely match).
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Resize the (now narrower) status columns ---
# Target stored column width is 13.4101845877511 characters; the COM
# ColumnWidth setter is quantized to the nearest 1/6 of a character by the
# host, so we request the width whose rounded result lands closest to the
# desired value (target minus the standard 5/6 character padding).
$newColumnWidth = 13.4101845877511 - (5/6)

$overview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C
